# Insert two new weekly price rows ("Cebollín" / Vega Modelo de Temuco) into
# the data table, pushing the existing rows 568-598 down to 570-600 and
# updating the used-range dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 568 (shifts 568..598 -> 570..600)
$ws.Range("568:569").Insert()

# --- New row 568 ---
$ws.Cells.Item(568, 1).Value  = 10
$ws.Cells.Item(568, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(568, 3).Value  = "La Araucanía"
$ws.Cells.Item(568, 4).Value  = 45041
$ws.Cells.Item(568, 5).Value  = 9
$ws.Cells.Item(568, 6).Value  = 100112037
$ws.Cells.Item(568, 7).Value  = "Cebollín"
$ws.Cells.Item(568, 8).Value  = "Sin especificar"
$ws.Cells.Item(568, 9).Value  = "Primera"
$ws.Cells.Item(568, 10).Value = 60
$ws.Cells.Item(568, 11).Value = 7000
$ws.Cells.Item(568, 12).Value = 7000
$ws.Cells.Item(568, 13).Value = 7000
$ws.Cells.Item(568, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(568, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(568, 16).Value = 583
$ws.Cells.Item(568, 17).Value = 12
$ws.Cells.Item(568, 18).Value = "Hortaliza"

# --- New row 569 ---
$ws.Cells.Item(569, 1).Value  = 10
$ws.Cells.Item(569, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(569, 3).Value  = "La Araucanía"
$ws.Cells.Item(569, 4).Value  = 45041
$ws.Cells.Item(569, 5).Value  = 9
$ws.Cells.Item(569, 6).Value  = 100112037
$ws.Cells.Item(569, 7).Value  = "Cebollín"
$ws.Cells.Item(569, 8).Value  = "Sin especificar"
$ws.Cells.Item(569, 9).Value  = "Primera"
$ws.Cells.Item(569, 10).Value = 90
$ws.Cells.Item(569, 11).Value = 7000
$ws.Cells.Item(569, 12).Value = 7000
$ws.Cells.Item(569, 13).Value = 7000
$ws.Cells.Item(569, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(569, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(569, 16).Value = 583
$ws.Cells.Item(569, 17).Value = 12
$ws.Cells.Item(569, 18).Value = "Hortaliza"
